# Applies the cryptos.xlsx data refresh described in the commit
# "Updated cryptos list on Thu Aug  1 06:44:27 UTC 2024 with GitHub Actions".
# Column D (Price) holds numeric-looking text (e.g. "64.414.63", "0.999") that
# must stay text, so its number format is forced to "@" before the value is
# written - otherwise the COM layer coerces it into a real number and drops
# formatting (e.g. "1.00" -> 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.414.63"
$ws.Range("E2").Value = "  -3.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.179.19"
$ws.Range("E3").Value = "  -4.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.10"
$ws.Range("E5").Value = "  -2.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.51"
$ws.Range("E6").Value = "  -7.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.611"
$ws.Range("E7").Value = "  -5.46%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.179.54"
$ws.Range("E9").Value = "  -4.12%  "
$ws.Range("E10").Value = "  -4.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.73"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.386"
$ws.Range("E12").Value = "  -4.30%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.724.78"
$ws.Range("E13").Value = "  -4.44%  "
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.425.45"
$ws.Range("E15").Value = "  -3.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.46"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("E17").Value = "  -3.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.180.81"
$ws.Range("E18").Value = "  -4.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "420.15"
$ws.Range("E19").Value = "  -2.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.82"
$ws.Range("E20").Value = "  -3.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.36"
$ws.Range("E21").Value = "  -3.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.04"
$ws.Range("E22").Value = "  -5.42%  "
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.22"
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("E25").Value = "  +2.13%  "
$ws.Range("E26").Value = "  -5.57%  "
$ws.Range("E27").Value = "  -7.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.72"
$ws.Range("E28").Value = "  -3.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "21.84"
$ws.Range("E30").Value = "  -2.78%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.83"
$ws.Range("E31").Value = "  -6.34%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.01"
$ws.Range("E33").Value = "  -4.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.33"
$ws.Range("E34").Value = "  -4.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.14"
$ws.Range("E35").Value = "  -5.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.59"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("E37").Value = "  -6.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.699.39"
$ws.Range("E38").Value = "  -6.40%  "
$ws.Range("E39").Value = "  -7.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.54"
$ws.Range("E40").Value = "  -8.63%  "
$ws.Range("E41").Value = "  -4.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.09"
$ws.Range("E42").Value = "  -2.97%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.711"
$ws.Range("E43").Value = "  -7.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.74"
$ws.Range("E44").Value = "  -4.97%  "
$ws.Range("E45").Value = "  -6.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0262"
$ws.Range("E46").Value = "  -3.83%  "
$ws.Range("E47").Value = "  -7.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "294.01"
$ws.Range("E48").Value = "  -7.31%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.02"
$ws.Range("E49").Value = "  -13.21%  "
$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0993"
$ws.Range("E51").Value = "  -4.86%  "
